$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the existing "hyperlink" cell style (used by column C) before we
# touch anything, so we can re-apply it after rebuilding the hyperlinks.
$hyperlinkStyle = $ws.Range("C2").Style
$normalStyle    = $ws.Range("A2").Style

# ---------------------------------------------------------------------------
# 1. Update the ID values for rows that were renamed (6,7,8,9)
#    (their e-mail addresses are rebuilt together with the hyperlinks below)
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "l123"
$ws.Range("B7").Value = "l456"
$ws.Range("B8").Value = "l789"
$ws.Range("B9").Value = "l135"

# ---------------------------------------------------------------------------
# 2. Row 14 name + GPA change (id/e-mail stay ll123 / ll123@qq.com)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Jelie"
$ws.Range("D14").Value = 2.98

# ---------------------------------------------------------------------------
# 3. Append the four brand-new rows (15-18)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Rumer"
$ws.Range("B15").Value = "ll456"
$ws.Range("D15").Value = 3.03

$ws.Range("A16").Value = "Kock"
$ws.Range("B16").Value = "ll789"
$ws.Range("D16").Value = 3.09

$ws.Range("A17").Value = "Naomy"
$ws.Range("B17").Value = "ll135"
$ws.Range("D17").Value = 3.11

$ws.Range("A18").Value = "Otis"
$ws.Range("B18").Value = "ll357"
$ws.Range("D18").Value = 3.14

for ($row = 15; $row -le 18; $row++) {
    $ws.Range("A$row").Style = $normalStyle
    $ws.Range("B$row").Style = $normalStyle
    $ws.Range("D$row").Style = $normalStyle
}

# ---------------------------------------------------------------------------
# 4. Rebuild the mailto hyperlinks for column C (rows 2-18) in row order, so
#    every link address matches the displayed e-mail text and relationship
#    ids stay sequential.
# ---------------------------------------------------------------------------
$ws.Range("C2").Hyperlinks.Delete()

$emails = @(
    "m129@qq.com", "m456@qq.com", "m789@qq.com", "m135@qq.com",
    "l123@qq.com", "l456@qq.com", "l789@qq.com", "l135@qq.com",
    "m123@qq.com", "m456@qq.com", "n789@qq.com", "n135@qq.com",
    "ll123@qq.com", "ll456@qq.com", "ll789@qq.com", "ll135@qq.com", "ll357@qq.com"
)

for ($i = 0; $i -lt $emails.Count; $i++) {
    $row = $i + 2
    $email = $emails[$i]
    $cell = $ws.Range("C$row")
    $cell.Value = $email
    $ws.Hyperlinks.Add($cell, "mailto:$email", [Type]::Missing, [Type]::Missing, $email) | Out-Null
    $cell.Style = $hyperlinkStyle
}

# ---------------------------------------------------------------------------
# 5. Misc view/selection tweak from the diff
# ---------------------------------------------------------------------------
$ws.Range("F13").Select()
